$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E both hold plain text in this sheet (prices like "1.796.33" and
# padded percentages like "  -1.22%  "). Assigning such strings via .Value can make
# Excel auto-detect a small subset of them (pure decimal-looking numbers, e.g.
# "309.20") as numeric. To guarantee they stay text -- matching the original
# inlineStr cells -- each touched D cell is temporarily switched to the "@" (text)
# number format before the assignment, then restored to the default "Normal" style
# afterwards so no stray formatting differences remain.

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.698.48"
$ws.Range("E2").Value = "  -1.15%  "
Set-TextValue $ws.Range("D3") "1.795.63"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue $ws.Range("D5") "309.08"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("E6").Value = "  -0.03%  "
Set-TextValue $ws.Range("D7") "0.4438"
$ws.Range("E7").Value = "  +5.42%  "
Set-TextValue $ws.Range("D8") "0.3679"
$ws.Range("E8").Value = "  +0.53%  "
Set-TextValue $ws.Range("D10") "0.8594"
$ws.Range("E10").Value = "  +2.21%  "
Set-TextValue $ws.Range("D11") "20.62"
$ws.Range("E11").Value = "  -0.73%  "
Set-TextValue $ws.Range("D12") "1.792.98"
$ws.Range("E12").Value = "  -1.28%  "
Set-TextValue $ws.Range("D13") "6.619"
$ws.Range("E13").Value = "  -0.02%  "
Set-TextValue $ws.Range("D14") "0.07063"
$ws.Range("E14").Value = "  +0.03%  "
Set-TextValue $ws.Range("D15") "91.55"
$ws.Range("E15").Value = "  +3.27%  "
Set-TextValue $ws.Range("D16") "5.267"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  -0.11%  "
Set-TextValue $ws.Range("D18") "0.000008684"
$ws.Range("E18").Value = "  -1.17%  "
Set-TextValue $ws.Range("D20") "14.79"
$ws.Range("E20").Value = "  -0.96%  "
Set-TextValue $ws.Range("D21") "26.715.38"
$ws.Range("E21").Value = "  -1.32%  "
Set-TextValue $ws.Range("D22") "5.163"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  +0.32%  "
Set-TextValue $ws.Range("D26") "2.176"
$ws.Range("E26").Value = "  -2.09%  "
Set-TextValue $ws.Range("D27") "18.40"
$ws.Range("E27").Value = "  +0.77%  "
Set-TextValue $ws.Range("D28") "5.181"
$ws.Range("E28").Value = "  -0.34%  "
Set-TextValue $ws.Range("D29") "117.13"
$ws.Range("E29").Value = "  +1.16%  "
Set-TextValue $ws.Range("D30") "0.08760"
$ws.Range("E30").Value = "  +0.01%  "
Set-TextValue $ws.Range("D31") "0.7389"
$ws.Range("E31").Value = "  +0.25%  "
Set-TextValue $ws.Range("D32") "1.151"
$ws.Range("E32").Value = "  -1.88%  "
Set-TextValue $ws.Range("D33") "4.440"
$ws.Range("E33").Value = "  +0.79%  "
Set-TextValue $ws.Range("D34") "2.904"
$ws.Range("E34").Value = "  -1.61%  "
Set-TextValue $ws.Range("D35") "1.000"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  -0.46%  "
Set-TextValue $ws.Range("D37") "0.01954"
$ws.Range("E37").Value = "  +0.10%  "
Set-TextValue $ws.Range("D38") "0.05183"
$ws.Range("E38").Value = "  -0.74%  "
Set-TextValue $ws.Range("D39") "0.5232"
$ws.Range("E39").Value = "  +4.15%  "
Set-TextValue $ws.Range("D40") "2.823"
$ws.Range("E40").Value = "  -1.70%  "
Set-TextValue $ws.Range("D41") "6.964"
$ws.Range("E41").Value = "  -3.95%  "
Set-TextValue $ws.Range("D42") "0.1682"
$ws.Range("E42").Value = "  -0.04%  "
Set-TextValue $ws.Range("D43") "0.5046"
$ws.Range("E43").Value = "  +6.79%  "
Set-TextValue $ws.Range("D44") "8.441"
$ws.Range("E44").Value = "  -1.19%  "
Set-TextValue $ws.Range("D45") "1.967"
$ws.Range("E45").Value = "  +5.30%  "
Set-TextValue $ws.Range("D46") "10.38"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("E48").Value = "  -0.07%  "
Set-TextValue $ws.Range("D49") "1.663"
$ws.Range("E49").Value = "  +1.34%  "
Set-TextValue $ws.Range("D50") "0.06287"
$ws.Range("E50").Value = "  -1.09%  "
Set-TextValue $ws.Range("D51") "0.9149"
$ws.Range("E51").Value = "  +1.61%  "
